$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "I" column (a spacer/rank column between "化学" and "生物") was removed
# entirely; everything to its right (生物/政治/地理) shifts one column left.
$ws.Columns("I").Delete()

# Re-point the conditional formatting range that used to cover A2:L88 so it
# matches the new (narrower) data extent A2:K88.
$fc = $ws.Range("A2:K88").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:K88"))

# Center the header row (row 1) horizontally & vertically.
$ws.Rows("1").HorizontalAlignment = -4108
$ws.Rows("1").VerticalAlignment = -4108

# Leave the same post-edit selection Excel would show after deleting the
# column: the column that slid into the old "J" position, fully selected.
$ws.Columns("J").EntireColumn.Select() | Out-Null
